# Split the run containing "...'Transaction ID' column and the merged dataset. "
# into two runs:
#   1) "...'Transaction ID' column " (keeps original run formatting)
#   2) "and the  "                    (new run, the word "merged dataset." is removed)
# on the "Approach" slide (slide 5), shape "Subtitle 2".

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text.IndexOf("and the merged dataset. ") -ge 0) {
                $targetSlide = $s
                $targetShape = $sh
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange
$tr.Replace("and the merged dataset. ", "and the  ", 1, 0, 0) | Out-Null
